# Re-order the "Recorded By" (column G) entries on the active sheet.
#
# Rule observed from the target diff:
#   - Split the cell text on ", ".
#   - If the literal entry "System" (exact case) is present, move it to the
#     front of the list, preserving the relative order of the remaining
#     entries.
#   - Otherwise (no "System" entry, but still multiple entries), reverse
#     the order of the entries.
#   - Single-entry cells are left untouched (nothing to reorder).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Text

    if ($value -eq $null -or $value -eq "") {
        continue
    }

    $parts = @($value -split ", ")

    if ($parts.Count -lt 2) {
        continue
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $newParts = @("System")
        foreach ($p in $parts) {
            if (-not $p.Equals("System")) {
                $newParts += $p
            }
        }
    } else {
        $newParts = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $newParts += $parts[$i]
        }
    }

    $newValue = $newParts -join ", "

    if (-not $newValue.Equals($value)) {
        $cell.Value = $newValue
    }
}
